# Project and Skill Data Changes to Support Bi-Directional Relationship
#
# Inserts a new "id" column at the front of the "project" sheet, populates
# it with the uuid for each project row, and adds a new skill uuid to the
# skills list of the most recent project (row 5 / Feb '21 entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- 1. Insert a new blank column before column A, shifting B:K (was A:J) right ---
$ws.Columns("A:A").Insert(-4161)

# Match the (approximate, best-fit) width Excel computed for the new "id" column
# once it is populated with full-length uuid strings.
$ws.Columns("A:A").ColumnWidth = 40.333333333333336

# --- 2. Populate the new "id" column ---
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = "259adedf-5510-42a8-9a28-c72754626b89"
$ws.Range("A3").Value = "4c5bdbad-b5d8-4bb2-bf49-f0244fd365f6"
$ws.Range("A4").Value = "2ed5500a-34ee-42a9-b053-77369cbc065a"
$ws.Range("A5").Value = "31793aed-8cbe-40e4-b0f6-1b90a098e999"

# Copy the existing header / data formatting from column B onto the new column A
# so the look (fonts, fills, alignment, borders) matches the rest of the table.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B2:B5").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3. Append a new skill id to the skills list of the Feb '21 project (row 5) ---
$newSkills = "[" + $nl + `
    '"e8a6997b-c4d1-4ffb-a214-4ed99487b858",' + $nl + `
    '"d0eb1478-f755-49be-97bd-1b98c27e4b4a",' + $nl + `
    '"0535d6cc-3815-46e7-bf33-406bb9e3b094",' + $nl + `
    '"115c0855-4c6e-4943-963a-7ff60139a932",' + $nl + `
    '"f56d0a66-2b46-4697-b53e-4aad36ea0fc5",' + $nl + `
    '"fc99bf5d-ead9-44eb-8de7-dd0add714411",' + $nl + `
    '"5d12e67a-f6db-48fd-848c-e267add62128",' + $nl + `
    '"f6a97e57-0b4e-459a-b8d7-0dffe189be5a",' + $nl + `
    '"79de5ad4-24ca-47c1-b18d-798c68705e1b",' + $nl + `
    '"48ec5901-f927-4b3d-8a12-4c0d67185367",' + $nl + `
    '"823b6fa1-5bcf-498e-bf5e-29387b6fac3b",' + $nl + `
    '"6a7cbf58-8987-44ab-8065-fc527ed632e2",' + $nl + `
    '"2017a078-6d1f-4705-af3a-f60e2289a8ee",' + $nl + `
    '"e2a283c2-85f3-4941-b38e-71cedb3e4dec",' + $nl + `
    '"3991bdd1-6fa1-4d3a-af05-919965dcaa3d"' + $nl + `
    "]"

$ws.Range("H5").Value = $newSkills

# --- 4. Update the active cell selection (shifted one column right, same as Excel would) ---
$ws.Range("H6").Select()

Write-Host "Edit complete"
